# Applies weekly update to the fruit/vegetable pricing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45134
$ws.Range("J2").Value = 50

# Row 3
$ws.Range("I3").Value = "Primera"
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2500
$ws.Range("P3").Value = 833

# Row 4
$ws.Range("D4").Value = 45149
$ws.Range("I4").Value = "Segunda"
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 2000
$ws.Range("P4").Value = 667

# Row 5
$ws.Range("D5").Value = 45145
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = 2500
$ws.Range("P5").Value = 833

# Row 6
$ws.Range("D6").Value = 45145
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("P6").Value = 667

# Row 7
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 1200
$ws.Range("L7").Value = 1300
$ws.Range("M7").Value = 1250
$ws.Range("P7").Value = 417

# Row 8
$ws.Range("D8").Value = 44838
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 1000
$ws.Range("P8").Value = 333

# Row 9
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 2500
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2500
$ws.Range("P9").Value = 833

# Row 10
$ws.Range("D10").Value = 45148
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 2000
$ws.Range("P10").Value = 667

# Row 11
$ws.Range("D11").Value = 45135
$ws.Range("J11").Value = 70
$ws.Range("K11").Value = 2500
$ws.Range("L11").Value = 2500
$ws.Range("M11").Value = 2500
$ws.Range("P11").Value = 833

# Row 12
$ws.Range("D12").Value = 44832
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 1200
$ws.Range("L12").Value = 1300
$ws.Range("M12").Value = 1250
$ws.Range("P12").Value = 417

# Row 13
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 1000
$ws.Range("P13").Value = 333

# Row 14
$ws.Range("D14").Value = 45146
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("P14").Value = 833

# Row 15
$ws.Range("D15").Value = 45146
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 2000
$ws.Range("P15").Value = 667

# Row 16
$ws.Range("D16").Value = 44846
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 1200
$ws.Range("L16").Value = 1300
$ws.Range("M16").Value = 1250
$ws.Range("P16").Value = 417

# Row 17
$ws.Range("D17").Value = 44846
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 1000
$ws.Range("M17").Value = 1000
$ws.Range("P17").Value = 333

# Row 18
$ws.Range("D18").Value = 45133
$ws.Range("J18").Value = 80
